$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.587.60"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").Value = "3.039.27"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.69"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.34"
$ws.Range("E6").Value = "  -2.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "3.037.48"
$ws.Range("E8").Value = "  -1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +4.21%  "

$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("E11").Value = "  -13.24%  "

$ws.Range("E12").Value = "  +3.73%  "

$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.36"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "3.537.14"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "63.652.38"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").Value = "3.044.33"
$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.49"
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("E23").Value = "  +8.40%  "

$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.43"
$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.03"
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.91"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0406"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "438.41"
$ws.Range("E38").Value = "  -6.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0810"
$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").Value = "2.993.84"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  -5.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.51"
$ws.Range("E45").Value = "  -3.02%  "

$ws.Range("E46").Value = "  +5.35%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.90"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("E50").Value = "  -1.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").Value = "  -0.03%  "
